$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 empty columns starting at column C, pushing the existing
# "English"/"Math" columns (C,D) to the right by 4 (to G,H).
$ws.Range("C1:F1").EntireColumn.Insert()
